# Update version string from the old build to the new release build.
#
# Old: mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)
# New: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$aboutSheet = $wb.Worksheets.Item("About")

# A2: "Version: ..."
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation text
$aboutSheet.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Chaili Coal Mine, China, M1829, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# --- "Boundaries and methane sources" sheet ---
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Find the "build_version" column and update every data row's value.
$headerRange = $dataSheet.Range("A1:Z1")
$buildVersionCol = $null
for ($c = 1; $c -le $headerRange.Columns.Count; $c++) {
    if ($dataSheet.Cells.Item(1, $c).Value2 -eq "build_version") {
        $buildVersionCol = $c
        break
    }
}

$usedRange = $dataSheet.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellValue = $dataSheet.Cells.Item($r, $buildVersionCol).Value2
    if ($cellValue -eq $oldVersion) {
        $dataSheet.Cells.Item($r, $buildVersionCol).Value = $newVersion
    }
}
